$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spain (row 6) - new cases/deaths figures updated
$ws.Range("E6").Value = 17780
$ws.Range("G6").Value = 213
$ws.Range("H6").Value = 1044

# Row 7 now shows Alemania (was Iran) with new data, row 8 now shows Iran
# (previous Alemania/Iran data) - i.e. Alemania and Iran swap their
# positions in the country list, and Alemania receives updated figures.
$ws.Range("A7").Value = "Alemania"
$ws.Range("B7").Value = 19711
$ws.Range("C7").Value = 4391
$ws.Range("D7").Value = 180
$ws.Range("E7").Value = 19472
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 59

$ws.Range("A8").Value = "Iran"
$ws.Range("B8").Value = 19644
$ws.Range("C8").Value = 1237
$ws.Range("D8").Value = 6745
$ws.Range("E8").Value = 11466
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 149
$ws.Range("H8").Value = 1433

# Estados Unidos (row 9)
$ws.Range("B9").Value = 16594
$ws.Range("C9").Value = 2805
$ws.Range("E9").Value = 16244
$ws.Range("G9").Value = 18
$ws.Range("H9").Value = 225

# Suiza (row 12)
$ws.Range("B12").Value = 5369
$ws.Range("C12").Value = 1147
$ws.Range("E12").Value = 5298

# Reino Unido (row 13)
$ws.Range("B13").Value = 3983
$ws.Range("C13").Value = 714
$ws.Range("E13").Value = 3741
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = 177

# Noruega (row 17)
$ws.Range("B17").Value = 1921
$ws.Range("C17").Value = 131
$ws.Range("E17").Value = 1913

# Timestamp update
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 19:16"
